$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.404.36'
$ws.Range('E2').Value = '  +2.13%  '
$ws.Range('D3').Value = '1.840.52'
$ws.Range('E3').Value = '  +1.92%  '
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.36'
$ws.Range('E5').Value = '  +4.10%  '
$ws.Range('E6').Value = '  +2.81%  '
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.97'
$ws.Range('E8').Value = '  +11.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.313'
$ws.Range('E9').Value = '  +7.88%  '
$ws.Range('E10').Value = '  +3.50%  '
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('D12').Value = '2.109.26'
$ws.Range('E12').Value = '  +2.00%  '
$ws.Range('D13').Value = '1.863.20'
$ws.Range('E13').Value = '  +3.08%  '
$ws.Range('E14').Value = '  +3.48%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.75'
$ws.Range('E15').Value = '  +8.49%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.672'
$ws.Range('E16').Value = '  +5.81%  '
$ws.Range('D17').Value = '35.442.57'
$ws.Range('E17').Value = '  +2.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.56'
$ws.Range('E18').Value = '  +3.63%  '
$ws.Range('E19').Value = '  +4.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '242.01'
$ws.Range('E20').Value = '  +0.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.01'
$ws.Range('E21').Value = '  +8.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.64'
$ws.Range('E22').Value = '  +13.45%  '
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('E24').Value = '  +3.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.31'
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.90'
$ws.Range('E26').Value = '  +2.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.74'
$ws.Range('E27').Value = '  +1.60%  '
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.60'
$ws.Range('E29').Value = '  +30.76%  '
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('D31').Value = '3.342.39'
$ws.Range('E31').Value = '  +37.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0562'
$ws.Range('E32').Value = '  +9.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.96'
$ws.Range('E33').Value = '  +5.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.07'
$ws.Range('E34').Value = '  +5.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.82'
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '94.72'
$ws.Range('E36').Value = '  +14.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.689'
$ws.Range('E37').Value = '  +7.14%  '
$ws.Range('E38').Value = '  +5.77%  '
$ws.Range('E39').Value = '  +4.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '15.31'
$ws.Range('E40').Value = '  +4.56%  '
$ws.Range('D41').Value = '1.326.29'
$ws.Range('E41').Value = '  +1.22%  '
$ws.Range('E42').Value = '  +7.19%  '
$ws.Range('E43').Value = '  +2.13%  '
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.27'
$ws.Range('E47').Value = '  +9.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0513'
$ws.Range('E48').Value = '  -1.00%  '
$ws.Range('D49').Value = '2.018.81'
$ws.Range('E49').Value = '  +2.51%  '
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '101.76'
$ws.Range('E51').Value = '  -0.21%  '
